$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ32251069",
    "summ32817601",
    "summ33402063",
    "summ34002909",
    "summ34563328",
    "summ35207633",
    "summ35882013",
    "summ36457433",
    "summ37014408"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $sheet = $wb.Worksheets.Item($i + 1)
    $sheet.Name = $newNames[$i]
}
